# M13 over ISO-25010 is vervallen.
#
# Maatregel M13 "Het project gebruikt ISO-25010 voor de specificatie van
# productkwaliteitseisen" is vervallen (row 40 of the "Self-assessment
# checklist" sheet). Removing the row shifts every row below it up by one,
# but this COM-interop runtime does not automatically re-home the
# worksheet's Comments / Hyperlinks collections when a row is deleted, so
# we capture their (row, column, content) up front, clear them, delete the
# row, and re-create them at their new row numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Self-assessment checklist")

$targetRow = 40

# ---- 1. capture existing comments (row, col, text) ----------------------
$commentData = New-Object System.Collections.ArrayList
foreach ($cm in $ws.Comments) {
    $p = $cm.Parent
    $info = @{ Row = $p.Row; Col = $p.Column; Text = $cm.Text() }
    [void]$commentData.Add($info)
}

# ---- 2. capture existing hyperlinks (row, col, subaddress, address) -----
$linkData = New-Object System.Collections.ArrayList
foreach ($h in $ws.Hyperlinks) {
    $r = $h.Range
    $info = @{ Row = $r.Row; Col = $r.Column; SubAddress = $h.SubAddress; Address = $h.Address }
    [void]$linkData.Add($info)
}

# ---- 3. clear all comments & hyperlinks so stale refs don't linger ------
while ($ws.Comments.Count -gt 0) {
    $ws.Comments.Item(1).Delete()
}
$ws.Hyperlinks.Delete()

# ---- 4. delete the row for M13 (shifts everything below up by one) ------
$ws.Rows.Item($targetRow).EntireRow.Delete()

# ---- 5. re-create the comments, skipping the one that lived on the ------
#         deleted row, shifting everything below it up by one row --------
foreach ($d in $commentData) {
    if ($d.Row -eq $targetRow) {
        continue
    }
    $newRow = $d.Row
    if ($newRow -gt $targetRow) {
        $newRow = $newRow - 1
    }
    $cell = $ws.Cells.Item($newRow, $d.Col)
    $cell.AddComment($d.Text)
}

# ---- 6. re-create the hyperlinks, same shifting rule ---------------------
foreach ($d in $linkData) {
    if ($d.Row -eq $targetRow) {
        continue
    }
    $newRow = $d.Row
    if ($newRow -gt $targetRow) {
        $newRow = $newRow - 1
    }
    $cell = $ws.Cells.Item($newRow, $d.Col)
    $ws.Hyperlinks.Add($cell, $d.Address, $d.SubAddress)
}
